$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values such as
# "0.9979" or "27.951.40" are not auto-converted to numbers/dates,
# matching the original inline-string cell contents.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.951.40'
$ws.Range("D3").Value = '1.912.23'
$ws.Range("E3").Value = '  +0.07%  '
$ws.Range("D4").Value = '0.9979'
$ws.Range("E4").Value = '  -0.64%  '
$ws.Range("D5").Value = '313.28'
$ws.Range("D6").Value = '0.9982'
$ws.Range("E6").Value = '  -0.55%  '
$ws.Range("D7").Value = '0.5010'
$ws.Range("E7").Value = '  +3.94%  '
$ws.Range("D8").Value = '0.3819'
$ws.Range("E8").Value = '  +0.26%  '
$ws.Range("D9").Value = '0.07321'
$ws.Range("E9").Value = '  -0.49%  '
$ws.Range("D10").Value = '0.9127'
$ws.Range("E10").Value = '  -2.28%  '
$ws.Range("E11").Value = '  +2.16%  '
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").Value = '0.07671'
$ws.Range("E12").Value = '  -1.79%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.873.58'
$ws.Range("E13").Value = '  -1.90%  '
$ws.Range("D14").Value = '5.491'
$ws.Range("E14").Value = '  -0.28%  '
$ws.Range("D15").Value = '92.89'
$ws.Range("E15").Value = '  +0.97%  '
$ws.Range("D16").Value = '0.9986'
$ws.Range("E16").Value = '  -0.59%  '
$ws.Range("D17").Value = '0.000008751'
$ws.Range("E17").Value = '  -1.35%  '
$ws.Range("D18").Value = '0.9990'
$ws.Range("E18").Value = '  -0.52%  '
$ws.Range("D19").Value = '27.983.50'
$ws.Range("E19").Value = '  -0.26%  '
$ws.Range("D21").Value = '5.185'
$ws.Range("E21").Value = '  +0.29%  '
$ws.Range("D22").Value = '2.161.17'
$ws.Range("E22").Value = '  +1.28%  '
$ws.Range("D23").Value = '10.87'
$ws.Range("E23").Value = '  -0.38%  '
$ws.Range("D24").Value = '6.630'
$ws.Range("E24").Value = '  -0.11%  '
$ws.Range("D25").Value = '153.21'
$ws.Range("E25").Value = '  -2.64%  '
$ws.Range("D26").Value = '1.847'
$ws.Range("E26").Value = '  -3.43%  '
$ws.Range("D27").Value = '2.211'
$ws.Range("E27").Value = '  +3.46%  '
$ws.Range("D28").Value = '18.44'
$ws.Range("E28").Value = '  -0.33%  '
$ws.Range("D29").Value = '115.59'
$ws.Range("E29").Value = '  -1.40%  '
$ws.Range("D30").Value = '4.937'
$ws.Range("E30").Value = '  -0.78%  '
$ws.Range("E31").Value = '  +0.77%  '
$ws.Range("D32").Value = '3.211'
$ws.Range("E32").Value = '  -2.14%  '
$ws.Range("D33").Value = '4.848'
$ws.Range("E33").Value = '  +3.97%  '
$ws.Range("D34").Value = '1.242'
$ws.Range("E34").Value = '  -1.25%  '
$ws.Range("D35").Value = '0.7833'
$ws.Range("E35").Value = '  +1.04%  '
$ws.Range("B36").Value = 'RenderToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D36").Value = '2.626'
$ws.Range("E36").Value = '  +0.38%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = '0.02086'
$ws.Range("E37").Value = '  +1.83%  '
$ws.Range("E38").Value = '  +2.54%  '
$ws.Range("E39").Value = '  -1.16%  '
$ws.Range("D40").Value = '0.5559'
$ws.Range("E40").Value = '  +0.44%  '
$ws.Range("D41").Value = '0.05290'
$ws.Range("E41").Value = '  -0.11%  '
$ws.Range("D42").Value = '6.898'
$ws.Range("E42").Value = '  -1.86%  '
$ws.Range("D43").Value = '113.90'
$ws.Range("E43").Value = '  +4.96%  '
$ws.Range("D44").Value = '8.549'
$ws.Range("E44").Value = '  +0.36%  '
$ws.Range("D45").Value = '0.1520'
$ws.Range("E45").Value = '  -0.62%  '
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = '0.4837'
$ws.Range("E46").Value = '  +0.06%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '10.59'
$ws.Range("E47").Value = '  -1.36%  '
$ws.Range("D48").Value = '0.9978'
$ws.Range("E48").Value = '  -0.60%  '
$ws.Range("E49").Value = '  -0.58%  '
$ws.Range("D50").Value = '67.75'
$ws.Range("E50").Value = '  -0.73%  '
$ws.Range("E51").Value = '  -0.39%  '

# Restore the default (Normal) style on column D now that the text
# values are safely stored, so no stray number-format style lingers.
$ws.Range("D2:D51").Style = "Normal"
